$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q3" worksheet right before the existing
#    "2022-Q2" worksheet, so the tab order becomes: 总计, 2022-Q3, 2022-Q2.
# ---------------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($oldQ2)
$q3.Name = "2022-Q3"

$zj = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the Q3 fund-holdings table.
#    Reuse the existing header / index-column formatting from "总计" so the
#    new sheet's styles line up with the rest of the workbook instead of
#    creating brand-new style records.
# ---------------------------------------------------------------------------
$zj.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$zj.Range("A2").Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Force the fund detail columns to be entered as text (so codes like
# "009258" keep their leading zero and decimals like "3.97" aren't rounded
# to a number), then drop back to the Normal style so no stray numeric
# format sticks around on the saved cells.
$q3.Range("B2:G5").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "673060"
$q3.Range("C2").Value = "西部利得景瑞灵活配置混合A"
$q3.Range("D2").Value = "3.97"
$q3.Range("E2").Value = "93.10"
$q3.Range("F2").Value = "5.54"
$q3.Range("G2").Value = "0.2199"
$q3.Range("H2").Value = 1

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "009258"
$q3.Range("C3").Value = "西部利得景瑞灵活配置混合C"
$q3.Range("D3").Value = "1.32"
$q3.Range("E3").Value = "93.10"
$q3.Range("F3").Value = "5.54"
$q3.Range("G3").Value = "0.0731"
$q3.Range("H3").Value = 1

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "011351"
$q3.Range("C4").Value = "金鹰年年邮益一年持有期混合A"
$q3.Range("D4").Value = "3.43"
$q3.Range("E4").Value = "34.33"
$q3.Range("F4").Value = "0.71"
$q3.Range("G4").Value = "0.0244"
$q3.Range("H4").Value = 9

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "011352"
$q3.Range("C5").Value = "金鹰年年邮益一年持有期混合C"
$q3.Range("D5").Value = "0.27"
$q3.Range("E5").Value = "34.33"
$q3.Range("F5").Value = "0.71"
$q3.Range("G5").Value = "0.0019"
$q3.Range("H5").Value = 9

$q3.Range("B2:G5").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: push the existing 2022-Q2 summary
#    row down to row 3 (bumping its index from 0 to 1) and write the new
#    2022-Q3 summary row into row 2.
# ---------------------------------------------------------------------------
$oldDate = $zj.Range("B2").Value2
$oldCount = $zj.Range("C2").Value2
$oldValue = $zj.Range("D2").Value2

$zj.Range("A2").Copy()
$zj.Range("A3").PasteSpecial(-4122)      # xlPasteFormats

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = $oldDate
$zj.Range("C3").Value = $oldCount
$zj.Range("D3").Value = $oldValue

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 4
$zj.Range("D2").Value = 0.32

# ---------------------------------------------------------------------------
# 4. Restore "2022-Q2" as the selected tab (it was the active sheet before
#    this edit).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Activate()
